$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "22.60", "1.001").
# Force text format first so Excel does not silently convert them to numbers,
# matching the workbook author's original inline-string price formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.253.50'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '1.792.61'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '325.69'
$ws.Range("E5").Value = '  -2.87%  '
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.4451'
$ws.Range("E7").Value = '  +13.41%  '
$ws.Range("D8").Value = '0.3729'
$ws.Range("E8").Value = '  +9.82%  '
$ws.Range("D9").Value = '44.52'
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("D10").Value = '1.145'
$ws.Range("E10").Value = '  +1.80%  '
$ws.Range("D11").Value = '0.07508'
$ws.Range("E11").Value = '  +3.75%  '
$ws.Range("D12").Value = '22.60'
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = '6.275'
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '7.525'
$ws.Range("E15").Value = '  +5.71%  '
$ws.Range("D16").Value = '1.788.67'
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = '0.00001087'
$ws.Range("E17").Value = '  +2.33%  '
$ws.Range("D18").Value = '0.06742'
$ws.Range("E18").Value = '  +1.93%  '
$ws.Range("D19").Value = '80.93'
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '17.52'
$ws.Range("E21").Value = '  +3.08%  '
$ws.Range("D22").Value = '6.317'
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("D23").Value = '28.243.11'
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = '11.73'
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").Value = '2.424'
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").Value = '20.39'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").Value = '152.00'
$ws.Range("E27").Value = '  -1.72%  '
$ws.Range("D28").Value = '2.355'
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").Value = '1.992.01'
$ws.Range("E29").Value = '  +2.49%  '
$ws.Range("D30").Value = '132.62'
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("D31").Value = '1.224'
$ws.Range("E31").Value = '  -4.44%  '
$ws.Range("D32").Value = '4.028'
$ws.Range("E32").Value = '  -1.22%  '
$ws.Range("D33").Value = '5.799'
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("D34").Value = '0.09392'
$ws.Range("E34").Value = '  +7.51%  '
$ws.Range("E35").Value = '  +9.80%  '
$ws.Range("D36").Value = '12.07'
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").Value = '0.06337'
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("D38").Value = '0.02328'
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("D39").Value = '5.157'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '0.6526'
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").Value = '8.289'
$ws.Range("E41").Value = '  +4.70%  '
$ws.Range("D42").Value = '1.470'
$ws.Range("E42").Value = '  -1.69%  '
$ws.Range("D43").Value = '1.209'
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = '14.01'
$ws.Range("E45").Value = '  +1.24%  '
$ws.Range("D46").Value = '0.6075'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").Value = '3.778'
$ws.Range("E47").Value = '  -1.25%  '
$ws.Range("D48").Value = '129.95'
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").Value = '2.017'
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("D50").Value = '0.07117'
$ws.Range("E50").Value = '  +1.40%  '
$ws.Range("D51").Value = '1.157'
$ws.Range("E51").Value = '  -0.23%  '

# Restore the default (Normal) style on column D so no stray text-format
# styling is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"

